# Add a new "2022-Q1" sheet (fund-level holdings) right before the "总计" sheet,
# and update the "总计" (totals) sheet with a new leading row summarizing 2022-Q1.
#
# NOTE: worksheet references returned by Worksheets.Item(...) track a sheet's
# position (index), not its identity. Once sheets are inserted/copied, any
# previously captured reference may silently start pointing at a different
# sheet. To stay safe we always re-fetch the sheet we need *by name* right
# before using it, instead of reusing an older variable.

$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, [string]$text) {
    # Force the cell to keep a literal text value (e.g. "84.98") instead of
    # Excel auto-coercing it into a floating point number.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ------------------------------------------------------------------
# 1) Create the new "2022-Q1" worksheet by duplicating the "2021-Q4"
#    sheet (so header/row styling is carried over identically) and
#    placing the copy immediately before "总计".
# ------------------------------------------------------------------
$sheetQ4 = $wb.Worksheets.Item("2021-Q4")
$sheetTotal = $wb.Worksheets.Item("总计")

$sheetQ4.Copy($sheetTotal)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# The source sheet has 4 data rows; the new sheet only needs 2, so
# drop rows 4 and 5.
$newSheet.Rows.Item(4).Resize(2).Delete()

# Row 2: fund 011328
Set-TextValue $newSheet.Cells.Item(2, 4) "84.98"
Set-TextValue $newSheet.Cells.Item(2, 5) "92.34"
Set-TextValue $newSheet.Cells.Item(2, 6) "3.69"
Set-TextValue $newSheet.Cells.Item(2, 7) "3.1358"
$newSheet.Cells.Item(2, 8).Value = 10

# Row 3: fund 011329
Set-TextValue $newSheet.Cells.Item(3, 4) "55.30"
Set-TextValue $newSheet.Cells.Item(3, 5) "92.34"
Set-TextValue $newSheet.Cells.Item(3, 6) "3.69"
Set-TextValue $newSheet.Cells.Item(3, 7) "2.0406"
$newSheet.Cells.Item(3, 8).Value = 10

# ------------------------------------------------------------------
# 2) Insert a new summary row into "总计" for 2022-Q1, pushing the
#    existing 2021-Q4 / 2021-Q3 rows down by one.
#    Re-fetch "总计" by name: the sheet insertion above shifted its
#    position, so any earlier reference to it is no longer valid.
# ------------------------------------------------------------------
$sheetTotal = $wb.Worksheets.Item("总计")
$sheetTotal.Rows.Item(2).Insert()

$sheetTotal.Cells.Item(2, 1).Value = 0
$sheetTotal.Cells.Item(2, 2).Value = "2022-Q1"
$sheetTotal.Cells.Item(2, 3).Value = 2
$sheetTotal.Cells.Item(2, 4).Value = 5.18

# Renumber the index column for the rows that shifted down.
$sheetTotal.Cells.Item(3, 1).Value = 1
$sheetTotal.Cells.Item(4, 1).Value = 2

# The inserted row picked up stray formatting; restore it to match the
# rest of the table (A column styled like the other index cells, the
# remaining columns left unstyled).
$sheetTotal.Cells.Item(3, 1).Copy()
$sheetTotal.Cells.Item(2, 1).PasteSpecial(-4122)
$sheetTotal.Range("B2:D2").Style = "Normal"

# Restore the originally active sheet/tab so we don't leave an unrelated
# view-state change behind.
$wb.Worksheets.Item("2021-Q3").Activate()
